$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Target dataset: rows 16-43, columns C (N Doc Trabajador), D (Nombre Trabajador),
# E (Periodo Mora), F (Valor Mora), G (Salario Basico).
# The table is re-sorted by Periodo Mora (ascending) then by worker (original grouping order),
# and CATERINA DEL CARMEN PALLARES GARCES (45690192) Salario Basico is updated from 908526 to 877803.
$data = @(
    ,@('CC', '45690192', 'CATERINA DEL CARMEN PALLARES GARCES', '2109', 36341, 877803)
    ,@('CC', '1047464230', 'JESUS ALBERTO JUNCO RIZO', '2109', 36341, 908526)
    ,@('CC', '1001803763', 'JORGE ANDRES BELTRAN GONZALEZ', '2109', 55120, 1378000)
    ,@('CC', '23238652', 'DAMASA RUIZ FONSECA', '2109', 36341, 908526)
    ,@('CC', '45690192', 'CATERINA DEL CARMEN PALLARES GARCES', '2110', 36341, 877803)
    ,@('CC', '1047464230', 'JESUS ALBERTO JUNCO RIZO', '2110', 36341, 908526)
    ,@('CC', '1001803763', 'JORGE ANDRES BELTRAN GONZALEZ', '2110', 55120, 1378000)
    ,@('CC', '23238652', 'DAMASA RUIZ FONSECA', '2110', 36341, 908526)
    ,@('CC', '45690192', 'CATERINA DEL CARMEN PALLARES GARCES', '2111', 36341, 877803)
    ,@('CC', '1047464230', 'JESUS ALBERTO JUNCO RIZO', '2111', 36341, 908526)
    ,@('CC', '1001803763', 'JORGE ANDRES BELTRAN GONZALEZ', '2111', 55120, 1378000)
    ,@('CC', '23238652', 'DAMASA RUIZ FONSECA', '2111', 36341, 908526)
    ,@('CC', '45690192', 'CATERINA DEL CARMEN PALLARES GARCES', '2112', 36341, 877803)
    ,@('CC', '1047464230', 'JESUS ALBERTO JUNCO RIZO', '2112', 36341, 908526)
    ,@('CC', '1001803763', 'JORGE ANDRES BELTRAN GONZALEZ', '2112', 55120, 1378000)
    ,@('CC', '23238652', 'DAMASA RUIZ FONSECA', '2112', 36341, 908526)
    ,@('CC', '45690192', 'CATERINA DEL CARMEN PALLARES GARCES', '2201', 36341, 877803)
    ,@('CC', '1047464230', 'JESUS ALBERTO JUNCO RIZO', '2201', 36341, 908526)
    ,@('CC', '1001803763', 'JORGE ANDRES BELTRAN GONZALEZ', '2201', 55120, 1378000)
    ,@('CC', '23238652', 'DAMASA RUIZ FONSECA', '2201', 36341, 908526)
    ,@('CC', '45690192', 'CATERINA DEL CARMEN PALLARES GARCES', '2202', 36341, 877803)
    ,@('CC', '1047464230', 'JESUS ALBERTO JUNCO RIZO', '2202', 36341, 908526)
    ,@('CC', '1001803763', 'JORGE ANDRES BELTRAN GONZALEZ', '2202', 55120, 1378000)
    ,@('CC', '23238652', 'DAMASA RUIZ FONSECA', '2202', 36341, 908526)
    ,@('CC', '45690192', 'CATERINA DEL CARMEN PALLARES GARCES', '2203', 31601, 877803)
    ,@('CC', '1047464230', 'JESUS ALBERTO JUNCO RIZO', '2203', 32707, 908526)
    ,@('CC', '1001803763', 'JORGE ANDRES BELTRAN GONZALEZ', '2203', 49608, 1378000)
    ,@('CC', '23238652', 'DAMASA RUIZ FONSECA', '2203', 32707, 908526)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 16 + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 2).Value = $rec[0]   # B: Tipo Doc Trabajador
    $ws.Cells.Item($row, 3).Value = $rec[1]   # C: N Doc Trabajador
    $ws.Cells.Item($row, 4).Value = $rec[2]   # D: Nombre Trabajador
    $ws.Cells.Item($row, 5).Value = $rec[3]   # E: Periodo Mora
    $ws.Cells.Item($row, 6).Value = $rec[4]   # F: Valor Mora
    $ws.Cells.Item($row, 7).Value = $rec[5]   # G: Salario Basico
}

Write-Output "Rows 16-43 updated"
